$d = $word.ActiveDocument

$replacements = @(
    @("2025-11-21 Friday", "2025-11-22 Saturday"),
    @("896×9=", "221×8="),
    @("381×7=", "706×8="),
    @("563×4=", "851×9="),
    @("541×2=", "346×6="),
    @("867×2=", "738×5="),
    @("510×9=", "201×6="),
    @("754×3=", "957×6="),
    @("290×2=", "576×3="),
    @("286×4=", "958×2="),
    @("361×4=", "585×9="),
    @("228×5=", "193×8="),
    @("253×6=", "626×6="),
    @("211×5=", "280×5="),
    @("861×3=", "953×9="),
    @("589×5=", "172×2="),
    @("229×5=", "969×2="),
    @("436×4=", "584×2="),
    @("443×6=", "234×5="),
    @("829×4=", "169×7="),
    @("658×2=", "145×8="),
    @("883×5=", "855×4="),
    @("464×5=", "199×8="),
    @("110×6=", "853×7="),
    @("596×3=", "287×4="),
    @("897×9=", "329×6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
